$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the three discontinued parts (delete bottom-up so row numbers
#    for the not-yet-deleted rows stay valid).
#    Row 44 = LCD Connector
#    Row 39 = Temp Sensor Connector
#    Row 30 = Processor Power Switch
# ---------------------------------------------------------------------------
$ws.Rows(44).Delete()
$ws.Rows(39).Delete()
$ws.Rows(30).Delete()

# ---------------------------------------------------------------------------
# 2. Hyperlinks: this engine does not re-anchor existing hyperlink ranges
#    when rows above them are deleted, so clear everything and re-create the
#    full, correct set of hyperlinks against the final row positions.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C7"),  "http://www.mouser.com/ProductDetail/Maxim-Integrated/MAX3077EESA+/?qs=sGAEpiMZZMuobhpKLk3hh6ov3TfCBqZhbNybjDy0atQ%3d")
$ws.Hyperlinks.Add($ws.Range("C9"),  "http://www.newark.com/stmicroelectronics/ld1117s33ctr/ic-ldo-volt-reg-3-3v-0-8a-sot/dp/89K0626?CMP=AFC-OP")
$ws.Hyperlinks.Add($ws.Range("C12"), "http://www.adafruit.com/products/790")
$ws.Hyperlinks.Add($ws.Range("C33"), "http://www.mouser.com/ProductDetail/Diodes-Incorporated/1N4148WS-7-F/?qs=sGAEpiMZZMtoHjESLttvkr74rFM1mfYSUkQNUJ8i7JM%3d")
$ws.Hyperlinks.Add($ws.Range("C31"), "http://www.mouser.com/ProductDetail/OSRAM-Opto-Semiconductors/LH-N974-KN-1/?qs=sGAEpiMZZMt82OzCyDsLFAV097Vn80XJzM0DIFS2How%3d")
$ws.Hyperlinks.Add($ws.Range("A3"),  "http://www.ti.com/tool/ek-tm4c123gxl")
$ws.Hyperlinks.Add($ws.Range("C8"),  "http://www.newark.com/nxp/74hc4052d-653/ic-analog-mux-dmux-dual-4-x-1/dp/78R7402")
$ws.Hyperlinks.Add($ws.Range("C30"), "http://www.newark.com/wurth-elektronik/742792701/ferrite-bead-0-05ohm-500ma-0402/dp/78R5663")
$ws.Hyperlinks.Add($ws.Range("C29"), "https://www.sparkfun.com/products/9806")
$ws.Hyperlinks.Add($ws.Range("C34"), "http://www.digikey.com/product-detail/en/CDSW4148-G/641-1459-1-ND/3511544")
$ws.Hyperlinks.Add($ws.Range("C37"), "https://www.sparkfun.com/products/643")

# ---------------------------------------------------------------------------
# 3. Rename "ATMega Programming Headers" -> "ATMega Programming Male Headers"
#    (after the deletions above this row now sits at row 41)
# ---------------------------------------------------------------------------
$ws.Range("A41").Value = "ATMega Programming Male Headers"

# Make sure every hyperlinked cell carries the workbook's built-in Hyperlink
# cell style (matches how the pre-existing links in the sheet look).
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("C7").Style = "Hyperlink"
$ws.Range("C8").Style = "Hyperlink"
$ws.Range("C9").Style = "Hyperlink"
$ws.Range("C12").Style = "Hyperlink"
$ws.Range("C29").Style = "Hyperlink"
$ws.Range("C30").Style = "Hyperlink"
$ws.Range("C31").Style = "Hyperlink"
$ws.Range("C33").Style = "Hyperlink"
$ws.Range("C34").Style = "Hyperlink"
$ws.Range("C37").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 4. View state: select C44 and scroll so row 10 is at the top, matching the
#    author's final cursor position/viewport.
# ---------------------------------------------------------------------------
$ws.Range("C44").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
